# Applies the "feat: get ERC721 listing events" edit:
# For the relevant swap rows, populate column L ("txHash") with the
# on-chain transaction hash associated with that listing/sale event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baazaar Sales")

$txHashes = @{
    21 = "0xc001874e440d76f77297fa9b5fc3dbb9c0b947095de7ec19251180b9e187b45f"
    22 = "0xf028a95e6ca63d28a45fca7c293e15431e2cda03a696009f35d5fe3f99b3ec34"
    23 = "0x135b23240803be669a729043d04e2a8cbc8c8333f1d06dabcec296433a958b97"
    24 = "0xb4e680693a96abe28ce43ab3870d04603e7aaca1e0e3bc45f64a22db30a95f94"
    25 = "0x965a1c3653b64f1314ceb7b679d84c6c4abdb039f6412c622740ed4453daa350"
    26 = "0xa10aa475db232ba858b40fbf4bb112f159c3a987c82ce5dbc23375147b2f4970"
    27 = "0xcba5c4307728b97b52ccbb0202318b2ae3830bad56f9874c955927fa7d0e24e9"
    30 = "0xe46aef57ab1277742cdb4615b445f84daa3d9c70167b408421afd483bef6b332"
    31 = "0x2900626516b016991b364d256a73dd9eb06fcbc936322ccfee69f9cb02408933"
    32 = "0x252f9446583fe02212b888e28cb3ab45564296b837c0a485d629c825380299a8"
    33 = "0x75d18e16e30be2eac97853b25325577de49c3102e14a200501a454b371b1eced"
    37 = "0x833c30b44bf502ca6f8d1fe3656bee1baee94d00fabbce187d2385c5b5164f8a"
    43 = "0xb2951b4893486cdec69fac7d58463cbed0e259ce11dc63f3fa60de30c919995c"
    48 = "0xdd0ea5ce8ca55b47eaf710662ee13eb4c02117fd7a7dc86218681f567ce885e0"
    49 = "0x3022d7ba3378fcf3761e94ce5aa595327010f241902489e5e77d3b06fa7e7e64"
    50 = "0xa30e02309b9c7ecc01c9059aef8e849dcd6ee37cd7d42ae6571af399553efee7"
    52 = "0x3777721d0a0367cc4ef054fad4b504923472cd18f799e74c132bb2451aeed4a9"
    53 = "0x248ea178da7022cdcda30a4458a7a7508c7800ead3833bd6b6c0f759dc84167d"
    54 = "0x5a895ea59a5bff39a0d62f2fa02bd32db177562974c15493fa0186d88fb363c2"
    55 = "0x92e55023eb0ef69cc3fd43c3c669b786daa2f4e4cbdf0699ef3e508f4e455533"
    56 = "0x968cccb0cb33c16d88475de655d43af79a35b0cec88ac6ad820c44eb67d1015c"
    57 = "0xff11d39511df28f8cd74de6f63215282af0c8b6a1bcdb86d7bed437c62eca38a"
    58 = "0xc808dea1de1ca80a346113521e8eee23e1ee504f50c2ff4a3300ec4c70d4b6c2"
    59 = "0x6d8462c0d8d216414877b54e9ba7dad2f9519e18d8070b1b2c10b23d5e550dae"
    60 = "0xf8b500320e2f3c7cc56aeeac8a6e162c68f5f0434e190605e108c7e970028ae8"
    61 = "0x6eee60083c0e375a40bcb46aabe2027a93e2eac01d99a38d293872abd020f3e6"
    62 = "0x6fde78f16413e8d7164c99a5d4c00ffd5cb2410f4983fba696e3e3f9eb6d66e3"
    63 = "0x771ce0823d559f266e9d5bc590510b47da83232804f5d2d7404ba89394bb2a59"
    64 = "0x6b6b3eaa0dd39229f7add81aa1467b80453b9dcfbc17a9c3d541682950114d5c"
    65 = "0xb8fe37e67bbf0f8f52f83d149a7feb5d2dad8a555cef3d6950cb989f944e5cb4"
    66 = "0xe03333eb0e8844ce9997af941a1cdaadcc0404bea48ee94a6a5484b645b2628a"
    67 = "0x09a4561fbda51f6f1d26bbccbd667eb9237a8cf00dd3af7e9e3ee3615ae83d3c"
    68 = "0x1321559efdb06cf82ed098bf1847a73ca706c853e0e0afb05a44c7b81fac930a"
    69 = "0xb22c77d7241cf858bf09c09a3a4b269ec06660f0ee3a7c063be5118b63dfd71a"
    70 = "0x7de46b820ba288c188c5a419f1d761a9d1fdb5a9fc78a7d7419a5f2d5b9e76ca"
    74 = "0xcd320730312cee012376ec831ecfff9a14cdda6267756f9f15549c5514ac48ce"
    75 = "0xb88e49974790e4dd7d4fb5fd8c87d90f0361734c86b27b11a4a125370dea909e"
    76 = "0x2f8784f31432dab15d6051eb0f2a548e9c7789d2e1efd6ec3c567a1d22731978"
    77 = "0x4a5933eec7afca5f190ca463a5add8dfc4a7945b726b8662f210e18dc9b769c6"
    78 = "0x8d79f25da6c0ea5e26a8d326361ad13a63dbe140c8711aa12ccd38bcdfc4253e"
    79 = "0x74f60c474411fd9df085256deb3242a824f17f141b6341da08374142df669b60"
    80 = "0x648d1c9fbd59a3e0f1f5f22b8c693117d42bfeae6a4d90fd1bbda2e608e3715f"
    81 = "0xea2abe76de1a7eb5615377440678d88d271f53239802f480182a0857511fe876"
    82 = "0x9d30d137c8a641897bd9bddc142e135b790b4b7e4997dcff3f1224c2c0d7f4f3"
    83 = "0xd5694aeb24ca8a2a5d03fe24724f3c5edbe61ab319f1323b95c511b178bdfffd"
    88 = "0xe7e93f3bad3d4506e3d0d582891946a73e703a1a81b243bd50b66d7c355e605c"
    91 = "0x748000d006c49e4861cef2aff6c72145f141d6e77b797bb9f6587b4780a76019"
}

foreach ($row in $txHashes.Keys) {
    $ws.Cells.Item([int]$row, 12).Value = $txHashes[$row]
}
